# Auto-update draw results: append the 2025-11-02 Pick 4 draw as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 47
$rng = $ws.Range("A" + $newRow + ":E" + $newRow)

# The existing rows store every value (dates, numeric-looking phase/result
# codes, timestamps) as plain text, not as real numbers/dates. Force the new
# row to Text format before writing so Excel doesn't auto-convert strings
# like "2025-11-02" or "251102" into a date serial / number.
$rng.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2025-11-02"
$ws.Cells.Item($newRow, 2).Value = "Pick 4"
$ws.Cells.Item($newRow, 3).Value = "251102"
$ws.Cells.Item($newRow, 4).Value = "0-9-6-7"
$ws.Cells.Item($newRow, 5).Value = "2025-11-02T21:36:07.048+04:00"

# Restore the default "Normal" style so the new cells don't carry an
# explicit style index, matching the unstyled cells used elsewhere in the
# sheet (the Text content itself is preserved - only the display format
# reverts to General).
$rng.Style = "Normal"
